# Updates cryptos list cell values (Coin/Link/Price/Volume(1h)) to match
# the latest scrape. Mirrors the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value as literal text even when it looks like a number
# (e.g. "1.00", "19.96") so Excel does not silently coerce it to a numeric
# cell and strip the formatting the source data relies on. NumberFormat is
# forced to Text just for the assignment, then the cell style is restored to
# "Normal" so no stray formatting is left behind.
function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '51.506.79'
$ws.Cells.Item(2, 5).Value = '  -0.24%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '2.816.47'
$ws.Cells.Item(3, 5).Value = '  +1.68%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
# Row 5
Set-TextValue 5 4 '352.16'
$ws.Cells.Item(5, 5).Value = '  +5.68%  '
# Row 6
Set-TextValue 6 4 '112.96'
$ws.Cells.Item(6, 5).Value = '  -2.90%  '
# Row 7
Set-TextValue 7 4 '0.570'
$ws.Cells.Item(7, 5).Value = '  +5.72%  '
# Row 8
Set-TextValue 8 4 '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.14%  '
# Row 9
Set-TextValue 9 4 '0.593'
$ws.Cells.Item(9, 5).Value = '  +3.14%  '
# Row 10
Set-TextValue 10 4 '41.34'
$ws.Cells.Item(10, 5).Value = '  -1.45%  '
# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.66%  '
# Row 12
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 12 4 '0.131'
$ws.Cells.Item(12, 5).Value = '  +1.13%  '
# Row 13
$ws.Cells.Item(13, 2).Value = 'Chainlink'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 13 4 '19.96'
$ws.Cells.Item(13, 5).Value = '  -1.63%  '
# Row 14
Set-TextValue 14 4 '7.69'
$ws.Cells.Item(14, 5).Value = '  +0.30%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '3.267.68'
$ws.Cells.Item(15, 5).Value = '  +2.12%  '
# Row 16
$ws.Cells.Item(16, 4).Value = '2.825.13'
$ws.Cells.Item(16, 5).Value = '  +2.21%  '
# Row 17
Set-TextValue 17 4 '0.889'
$ws.Cells.Item(17, 5).Value = '  -0.35%  '
# Row 18
$ws.Cells.Item(18, 4).Value = '51.363.67'
$ws.Cells.Item(18, 5).Value = '  -0.46%  '
# Row 19
Set-TextValue 19 4 '7.37'
$ws.Cells.Item(19, 5).Value = '  +7.21%  '
# Row 20
$ws.Cells.Item(20, 5).Value = '  -4.52%  '
# Row 21
Set-TextValue 21 4 '13.37'
$ws.Cells.Item(21, 5).Value = '  -1.32%  '
# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0990'
$ws.Cells.Item(22, 5).Value = '  +1.16%  '
# Row 23
Set-TextValue 23 4 '270.35'
$ws.Cells.Item(23, 5).Value = '  -2.91%  '
# Row 24
Set-TextValue 24 4 '69.55'
$ws.Cells.Item(24, 5).Value = '  -0.34%  '
# Row 25
Set-TextValue 25 4 '2.75'
$ws.Cells.Item(25, 5).Value = '  +1.92%  '
# Row 26
Set-TextValue 26 4 '26.62'
$ws.Cells.Item(26, 5).Value = '  -0.85%  '
# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.03%  '
# Row 28
Set-TextValue 28 4 '10.28'
$ws.Cells.Item(28, 5).Value = '  +0.85%  '
# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.86%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  -2.57%  '
# Row 31
$ws.Cells.Item(31, 2).Value = 'OKB'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 31 4 '50.61'
$ws.Cells.Item(31, 5).Value = '  +0.74%  '
# Row 32
$ws.Cells.Item(32, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 32 4 '33.95'
$ws.Cells.Item(32, 5).Value = '  -3.48%  '
# Row 33
Set-TextValue 33 4 '5.80'
$ws.Cells.Item(33, 5).Value = '  +3.97%  '
# Row 34
Set-TextValue 34 4 '0.0441'
$ws.Cells.Item(34, 5).Value = '  +24.80%  '
# Row 35
Set-TextValue 35 4 '0.0821'
$ws.Cells.Item(35, 5).Value = '  -0.08%  '
# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.08%  '
# Row 37
Set-TextValue 37 4 '4.96'
$ws.Cells.Item(37, 5).Value = '  -1.13%  '
# Row 38
$ws.Cells.Item(38, 5).Value = '  -1.67%  '
# Row 39
Set-TextValue 39 4 '3.17'
$ws.Cells.Item(39, 5).Value = '  -2.33%  '
# Row 40
Set-TextValue 40 4 '18.00'
$ws.Cells.Item(40, 5).Value = '  -5.56%  '
# Row 41
Set-TextValue 41 4 '23.70'
$ws.Cells.Item(41, 5).Value = '  +2.09%  '
# Row 42
$ws.Cells.Item(42, 5).Value = '  +2.15%  '
# Row 43
Set-TextValue 43 4 '126.26'
$ws.Cells.Item(43, 5).Value = '  -0.90%  '
# Row 44
Set-TextValue 44 4 '2.50'
$ws.Cells.Item(44, 5).Value = '  +1.42%  '
# Row 45
Set-TextValue 45 4 '2.29'
$ws.Cells.Item(45, 5).Value = '  -0.51%  '
# Row 46
$ws.Cells.Item(46, 4).Value = '2.077.77'
$ws.Cells.Item(46, 5).Value = '  -0.64%  '
# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.53%  '
# Row 48
$ws.Cells.Item(48, 5).Value = '  +3.57%  '
# Row 49
Set-TextValue 49 4 '5.65'
$ws.Cells.Item(49, 5).Value = '  +1.77%  '
# Row 50
Set-TextValue 50 4 '0.917'
$ws.Cells.Item(50, 5).Value = '  +3.89%  '
# Row 51
Set-TextValue 51 4 '60.60'
$ws.Cells.Item(51, 5).Value = '  +0.31%  '
